$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2025-07-01 Tuesday" "2025-07-02 Wednesday"

Replace-Text "661÷7=" "535÷6="
Replace-Text "354÷2=" "387÷4="
Replace-Text "171÷7=" "406÷9="
Replace-Text "424÷4=" "370÷7="
Replace-Text "683÷7=" "895÷9="
Replace-Text "171÷8=" "850÷4="
Replace-Text "650÷2=" "376÷4="
Replace-Text "184÷3=" "546÷2="
Replace-Text "236÷6=" "568÷8="
Replace-Text "149÷3=" "385÷5="
Replace-Text "987÷4=" "537÷5="
Replace-Text "800÷7=" "791÷4="
Replace-Text "997÷9=" "521÷4="
Replace-Text "858÷2=" "161÷5="
Replace-Text "584÷7=" "698÷2="
Replace-Text "869÷8=" "103÷8="
Replace-Text "223÷6=" "771÷5="
Replace-Text "348÷8=" "441÷8="
Replace-Text "442÷4=" "180÷7="
Replace-Text "404÷9=" "480÷2="
Replace-Text "764÷2=" "438÷5="
Replace-Text "659÷7=" "272÷5="
Replace-Text "987÷8=" "973÷5="
Replace-Text "736÷9=" "445÷8="
Replace-Text "202÷3=" "214÷9="
